$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.620.60'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '2.306.64'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.60'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.623'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.59%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -1.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.72'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.60%  '
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.54'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.109'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.01'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.39'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('D16').Value = '2.656.48'
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('D17').Value = '2.306.16'
$ws.Range('E17').Value = '  -0.41%  '
$ws.Range('D18').Value = '42.701.89'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('B19').Value = 'InternetComputer(DFINITY)'
$ws.Range('C19').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.50'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +32.82%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.57'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.68%  '
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.18'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.92%  '
$ws.Range('E23').Value = '  -2.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.23%  '
$ws.Range('E25').Value = '  -2.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.01'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.97'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.35'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.82'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +14.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.69'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '37.56'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '165.88'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0886'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.60%  '
$ws.Range('E34').Value = '  -3.51%  '
$ws.Range('E35').Value = '  -1.15%  '
$ws.Range('E36').Value = '  -2.14%  '
$ws.Range('E37').Value = '  -1.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0355'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.72'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.71'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '70.54'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '95.54'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.56%  '
$ws.Range('E44').Value = '  +0.57%  '
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.38'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.66%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '115.79'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '81.21'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.29%  '
$ws.Range('D49').Value = '1.681.72'
$ws.Range('E49').Value = '  +4.01%  '
$ws.Range('E50').Value = '  -1.91%  '
$ws.Range('E51').Value = '  -2.06%  '

Write-Output "Applied cryptos update"